# AutoCommit_14 июня 2024 г. 9:27:02_SibNout2023
#
# Grade-sheet touch-up: a handful of students' homework/lab marks are
# filled in (several 2's become 5's, a couple of blank cells get a 5, a
# "Хочу 4" flag is set), and the colour-scale conditional formatting that
# already decorates column D is extended to the newly-edited row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 11 - Жуков Никита: ДЗ_1..ДЗ_3 raised 2 -> 5, ДЗ_4 cleared out,
# column "7" and the lab column filled in, "Хочу 4" flag added.
# ---------------------------------------------------------------------
$ws.Range("C11:E11").Value = 5
$ws.Range("F11").Clear()
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 5
$ws.Range("K11").Value = "Хочу 4"

# C11:E11 pick up the same highlighted (double thick border) look already
# used on row 14 - copy that formatting across.
$ws.Range("C14:E14").Copy() | Out-Null
$ws.Range("C11:E11").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Row 13 - Иванова Снежана: column "7" + both lab columns filled in.
# ---------------------------------------------------------------------
$ws.Range("H13").Value = 5
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 5
$ws.Range("I29:J29").Copy() | Out-Null
$ws.Range("I13:J13").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Row 15 - Косова Мария: columns "6"/"7" + both lab columns filled in.
# ---------------------------------------------------------------------
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 5
$ws.Range("J15").Value = 5
$ws.Range("I29:J29").Copy() | Out-Null
$ws.Range("I15:J15").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Row 17 - Масленникова Анастасия: column "7" + first lab column filled in.
# ---------------------------------------------------------------------
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 5
$ws.Range("I24").Copy() | Out-Null
$ws.Range("I17").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Conditional formatting: add the same red/yellow/green colour-scale rule
# used on D4/D14/J4:J31 to D11 as well.
# ---------------------------------------------------------------------
$cf = $ws.Range("D11").FormatConditions.AddColorScale(3)
$cf.SetFirstPriority()

# ---------------------------------------------------------------------
# Leave the cursor on the last-touched cell, like the author did.
# ---------------------------------------------------------------------
$ws.Range("K15").Select()
